$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: column B ("Contraseña") stores the DPI digits as TEXT. A plain
# numeric-looking string gets auto-converted to a number by Excel, so a
# leading apostrophe is used to force text entry; the resulting "quote
# prefix" cell formatting is cleared right away so no stray formatting is
# left behind on the cell.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# Row 2 - Jose Gonzalez: DPI/password (text) + DPI (numeric) + phone number updated
Set-TextValue $ws.Range("B2") "9375847119455"
$ws.Range("G2").Value = 9375847119455
$ws.Range("I2").Value = 22211124

# Row 3 - Patricia Medina: DPI/password (text) + DPI (numeric) + phone number updated
Set-TextValue $ws.Range("B3") "6630177783950"
$ws.Range("G3").Value = 6630177783950
$ws.Range("I3").Value = 99901235

# Row 4 - Rosa Rivera: DPI/password (text) + DPI (numeric) + phone number updated
Set-TextValue $ws.Range("B4") "2984797705615"
$ws.Range("G4").Value = 2984797705615
$ws.Range("I4").Value = 10123478

# Row 5 - Julia Ruiz: DPI/password (text) + DPI (numeric) + phone number updated
Set-TextValue $ws.Range("B5") "6991634010529"
$ws.Range("G5").Value = 6991634010529
$ws.Range("I5").Value = 21234589

# Row 6 - Carolina Castro: DPI/password (text) + DPI (numeric) + phone number + email updated
Set-TextValue $ws.Range("B6") "2270770502250"
$ws.Range("G6").Value = 2270770502250
$ws.Range("I6").Value = 32345690
$ws.Range("J6").Value = "12@gmail.com"

# Row 7 - Miguel Garcia: entire row removed
$ws.Rows(7).Delete()

$ws.Range("H17").Select()
